$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Summary": row 2, columns B..L
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6226591760299626
$wsSummary.Range("C2").Value = 0.5746864310148233
$wsSummary.Range("D2").Value = 0.9438202247191011
$wsSummary.Range("E2").Value = 0.7143869596031184
$wsSummary.Range("F2").Value = 0.8363757052771325
$wsSummary.Range("G2").Value = 0.921065579531876
$wsSummary.Range("H2").Value = 0.7968410273674761
$wsSummary.Range("I2").Value = 504
$wsSummary.Range("J2").Value = 373
$wsSummary.Range("K2").Value = 161
$wsSummary.Range("L2").Value = 30

# ---------------------------------------------------------------
# Sheet "Classification Report": rows 2..6, columns B..D (E unchanged)
# ---------------------------------------------------------------
$wsClass = $wb.Worksheets.Item("Classification Report")

$wsClass.Range("B2").Value = 0.8429319371727748
$wsClass.Range("C2").Value = 0.301498127340824
$wsClass.Range("D2").Value = 0.4441379310344827

$wsClass.Range("B3").Value = 0.5746864310148233
$wsClass.Range("C3").Value = 0.9438202247191011
$wsClass.Range("D3").Value = 0.7143869596031184

$wsClass.Range("B4").Value = 0.6226591760299626
$wsClass.Range("C4").Value = 0.6226591760299626
$wsClass.Range("D4").Value = 0.6226591760299626
$wsClass.Range("E4").Value = 0.6226591760299626

$wsClass.Range("B5").Value = 0.7088091840937991
$wsClass.Range("C5").Value = 0.6226591760299626
$wsClass.Range("D5").Value = 0.5792624453188006

$wsClass.Range("B6").Value = 0.7088091840937991
$wsClass.Range("C6").Value = 0.6226591760299626
$wsClass.Range("D6").Value = 0.5792624453188006

# ---------------------------------------------------------------
# Sheet "Confusion Matrix": rows 2..3, columns B..C
# ---------------------------------------------------------------
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

$wsConf.Range("B2").Value = 161
$wsConf.Range("C2").Value = 373

$wsConf.Range("B3").Value = 30
$wsConf.Range("C3").Value = 504
